$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 3336
$ws.Range("C4").Value = 575
$ws.Range("D4").Value = 2652
$ws.Range("E4").Value = 3145
$ws.Range("F4").Value = 4463
$ws.Range("G4").Value = 93
$ws.Range("H4").Value = 32
$ws.Range("I4").Value = 57
$ws.Range("J4").Value = 84
$ws.Range("K4").Value = 159
$ws.Range("L4").Value = 9698
$ws.Range("M4").Value = 1252
$ws.Range("N4").Value = 7868
$ws.Range("O4").Value = 9528
$ws.Range("P4").Value = 12385

$ws.Range("B5").Value = 2406
$ws.Range("C5").Value = 226
$ws.Range("D5").Value = 2183
$ws.Range("E5").Value = 2371
$ws.Range("F5").Value = 2897
$ws.Range("G5").Value = 57
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 55
$ws.Range("J5").Value = 56
$ws.Range("K5").Value = 65
$ws.Range("L5").Value = 12016
$ws.Range("M5").Value = 1548
$ws.Range("N5").Value = 10214
$ws.Range("O5").Value = 11128
$ws.Range("P5").Value = 14543

$ws.Range("B6").Value = 45
$ws.Range("C6").Value = 37
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 61
$ws.Range("F6").Value = 116
$ws.Range("G6").Value = 156
$ws.Range("H6").Value = 12
$ws.Range("I6").Value = 148
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 187
$ws.Range("L6").Value = 4298
$ws.Range("M6").Value = 834
$ws.Range("N6").Value = 3665
$ws.Range("O6").Value = 4041
$ws.Range("P6").Value = 6728

$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 9
$ws.Range("G7").Value = 445
$ws.Range("H7").Value = 32
$ws.Range("I7").Value = 416
$ws.Range("J7").Value = 427
$ws.Range("K7").Value = 501
$ws.Range("L7").Value = 6302
$ws.Range("M7").Value = 244
$ws.Range("N7").Value = 6055
$ws.Range("O7").Value = 6183
$ws.Range("P7").Value = 6852

$ws.Range("B8").Value = 2303
$ws.Range("C8").Value = 139
$ws.Range("D8").Value = 2110
$ws.Range("E8").Value = 2331
$ws.Range("F8").Value = 2485
$ws.Range("G8").Value = 394
$ws.Range("H8").Value = 33
$ws.Range("I8").Value = 366
$ws.Range("J8").Value = 377
$ws.Range("K8").Value = 458
$ws.Range("L8").Value = 10080
$ws.Range("M8").Value = 1174
$ws.Range("N8").Value = 8642
$ws.Range("O8").Value = 9976
$ws.Range("P8").Value = 12438

$ws.Range("B9").Value = 3097
$ws.Range("C9").Value = 519
$ws.Range("D9").Value = 2615
$ws.Range("E9").Value = 2872
$ws.Range("F9").Value = 4279
$ws.Range("G9").Value = 807
$ws.Range("H9").Value = 42
$ws.Range("I9").Value = 756
$ws.Range("J9").Value = 807
$ws.Range("K9").Value = 898
$ws.Range("L9").Value = 19007
$ws.Range("M9").Value = 4162
$ws.Range("N9").Value = 15208
$ws.Range("O9").Value = 17189
$ws.Range("P9").Value = 26782

$ws.Range("B10").Value = 2178
$ws.Range("C10").Value = 1515
$ws.Range("D10").Value = 124
$ws.Range("E10").Value = 2138
$ws.Range("F10").Value = 5161
$ws.Range("G10").Value = 875
$ws.Range("H10").Value = 53
$ws.Range("I10").Value = 812
$ws.Range("J10").Value = 856
$ws.Range("K10").Value = 997
$ws.Range("L10").Value = 18687
$ws.Range("M10").Value = 3388
$ws.Range("N10").Value = 12177
$ws.Range("O10").Value = 19438
$ws.Range("P10").Value = 23026

$ws.Range("B11").Value = 75
$ws.Range("C11").Value = 19
$ws.Range("D11").Value = 62
$ws.Range("E11").Value = 68
$ws.Range("F11").Value = 131
$ws.Range("G11").Value = 47
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 46
$ws.Range("J11").Value = 47
$ws.Range("K11").Value = 49
$ws.Range("L11").Value = 3166
$ws.Range("M11").Value = 181
$ws.Range("N11").Value = 2885
$ws.Range("O11").Value = 3158
$ws.Range("P11").Value = 3558

$ws.Range("B12").Value = 2292
$ws.Range("C12").Value = 229
$ws.Range("D12").Value = 2144
$ws.Range("E12").Value = 2213
$ws.Range("F12").Value = 2957
$ws.Range("G12").Value = 16
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 15
$ws.Range("J12").Value = 16
$ws.Range("K12").Value = 16
$ws.Range("L12").Value = 6333
$ws.Range("M12").Value = 722
$ws.Range("N12").Value = 5234
$ws.Range("O12").Value = 6345
$ws.Range("P12").Value = 7363

$ws.Range("B13").Value = 3123
$ws.Range("C13").Value = 524
$ws.Range("D13").Value = 2391
$ws.Range("E13").Value = 3113
$ws.Range("F13").Value = 3789
$ws.Range("G13").Value = 1110
$ws.Range("H13").Value = 45
$ws.Range("I13").Value = 1053
$ws.Range("J13").Value = 1109
$ws.Range("K13").Value = 1227
$ws.Range("L13").Value = 19977
$ws.Range("M13").Value = 2903
$ws.Range("N13").Value = 15693
$ws.Range("O13").Value = 20432
$ws.Range("P13").Value = 23634
